$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same "想去人数" (column F) data
# that was refreshed in this commit.
$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    4  = 1643
    5  = 633
    6  = 1103
    8  = 11602
    12 = 376
    14 = 815
    15 = 12416
    16 = 13171
    21 = 251
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
